$wb = $excel.ActiveWorkbook

# --- 1. Insert a new "State" column into hotel_info, between Hotel_Name and City ---
$hotel = $wb.Worksheets.Item("hotel_info")
$hotel.Columns("C:C").Insert()
$hotel.Range("C1").Value2 = "State"
$hotel.Range("C2").Value2 = "Louisiana"

# --- 2. Reorder sheet tabs: review_info moves before hotel_info ---
$review = $wb.Worksheets.Item("review_info")
$review.Move($hotel)
